# Insert a new weekly price record as row 532, shifting all existing
# records from row 532 onward down by one row (dimension grows from
# A1:R604 to A1:R605).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push the existing rows 532..604 down to 533..605
$ws.Rows.Item(532).Insert()

# Populate the newly inserted row 532 with the new record
$ws.Range("A532").Value = 3
$ws.Range("B532").Value = "Femacal de La Calera"
$ws.Range("C532").Value = "Coquimbo"
$ws.Range("D532").Value = 45127
$ws.Range("E532").Value = 5
$ws.Range("F532").Value = 100112043
$ws.Range("G532").Value = "Pepino ensalada"
$ws.Range("H532").Value = "Sin especificar"
$ws.Range("I532").Value = "Primera"
$ws.Range("J532").Value = 90
$ws.Range("K532").Value = 10500
$ws.Range("L532").Value = 11000
$ws.Range("M532").Value = 10778
$ws.Range("N532").Value = "$/caja 60 unidades"
$ws.Range("O532").Value = "Región de Arica y Parinacota"
$ws.Range("P532").Value = 180
$ws.Range("Q532").Value = 60
$ws.Range("R532").Value = "Hortaliza"
